# Insert a new weekly price record for Zanahoria (Terminal Hortofrutícola Agro
# Chillán) above the current row 219. Excel shifts every row from 219..271
# down to 220..272, preserving their existing contents, and grows the used
# range from A1:R271 to A1:R272.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(219).Insert()

$ws.Cells.Item(219, 1).Value = 7
$ws.Cells.Item(219, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(219, 3).Value = "Ñuble"
$ws.Cells.Item(219, 4).Value = 44641
$ws.Cells.Item(219, 5).Value = 16
$ws.Cells.Item(219, 6).Value = 100114013
$ws.Cells.Item(219, 7).Value = "Zanahoria"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 120
$ws.Cells.Item(219, 11).Value = 7000
$ws.Cells.Item(219, 12).Value = 7500
$ws.Cells.Item(219, 13).Value = 7250
$ws.Cells.Item(219, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(219, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(219, 16).Value = 362
$ws.Cells.Item(219, 17).Value = 20
$ws.Cells.Item(219, 18).Value = "Hortaliza"
